# Add the "Cell", "Email", "Parties" column headings next to the existing
# "Name" heading in row 1, and a sample data row underneath (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Row 1: column headings --------------------------------------------
$ws.Range("B1").Value = "Cell"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Parties"

# Give the new "Email"/"Parties" headings the same look as the existing
# "Name" heading (copy its formatting across). "Cell" (B1) keeps the
# plain/default look, matching the sample row below.
$ws.Range("A1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# -- Row 2: example person data -----------------------------------------
$ws.Range("A2").Value = "Example person"
$ws.Range("B2").Value = "072 042 2222"
$ws.Range("C2").Value = "test@example.com"
$ws.Range("D2").Value = "EFF, DA"

# Store the row as text so values such as the phone number / party list
# are preserved verbatim (not reinterpreted as numbers/dates).
$ws.Range("A2:D2").NumberFormat = "@"
